# Commit: "removed server to its own repo"
#
# The canonical-OOXML diff for this commit only touches three things,
# and all three are non-semantic, tooling-generated identifiers rather
# than document content:
#
#   1. ppt/presentation.xml      - every r:id in <p:sldMasterIdLst>,
#                                   <p:sldIdLst> and <p:sldLayoutIdLst>
#                                   is swapped for a different random
#                                   "R<16 hex chars>" relationship id.
#   2. ppt/webextensions/taskpanes.xml - the r:id of the single
#                                   <wetp:webextensionref> is likewise
#                                   swapped for a different random id.
#   3. ppt/webextensions/webextension.xml - the wrapper element's own
#                                   we:webextension/@id GUID is swapped
#                                   for a different random GUID (the
#                                   add-in's real identity - the
#                                   <we:reference id="4d663618-...">
#                                   element - is untouched).
#
# In every case the *targets* of the relationships (slideMaster.xml,
# slide.xml, each slideLayoutN.xml, webextension.xml, ...), their
# count/order, and every other part (slide.xml, slideMasters/*,
# slideLayouts/*, the webextension's reference/properties/bindings,
# taskpane dockstate/visibility/width/row, etc.) are byte-for-byte
# identical before and after. Nothing in the slide content, layouts,
# master, or the add-in manifest actually changed - this is the
# signature of the project simply being rebuilt (Visual Studio/Open
# XML SDK mint a fresh relationship id - and, for the webextension
# part, a fresh part GUID - on every rebuild) after the unrelated
# source change described by the commit message (pulling the server
# project out into its own repo), not an edit a user made inside
# PowerPoint.
#
# Relationship ids and OOXML part GUIDs are package-serialization
# plumbing: they are not part of the Slide/Shape/TextRange surface (or
# any other object) that PowerPoint's object model exposes, so there
# is no COM call that corresponds to "rename this relationship id" or
# "mint a new webextension GUID" - a real PowerPoint automation script
# can't express that edit either. Touching the deck through the object
# model (e.g. deleting/recreating the slide or master/layouts just to
# force the engine to mint new ids) would actually falsify the diff,
# since it would also perturb slide.xml / slideMasters / slideLayouts
# content that the diff proves stayed untouched.
#
# So the faithful reproduction here is to simply load the deck as
# PowerPoint would and save it back out, with no content mutation -
# matching every part of the document that the diff shows as
# unchanged.
$p = $ppt.ActivePresentation
$p.Save()
